# Project "Sample Project" save: Rules sheet, cell B11 changes from the
# shared string "R40" to the (new) shared string "1".
#
# A plain  $ws.Range("B11").Value = "1"  would let Excel's usual
# numeric-literal auto-detection turn "1" into the *number* 1 (stored as
# <v>1</v> with no t="s"), which is not what the source workbook shows
# (B11 keeps t="s", i.e. it is still a text/string cell referencing a
# shared-string table entry, just a different one). To force the literal
# text "1" to land as a string value, build it via a text formula
# ("1") in a scratch cell, then paste-special just the resulting value
# into B11 - this preserves the cell's existing style/format and stores
# a genuine string "1" rather than a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)   # xlPasteValues
$excel.CutCopyMode = 0
$scratch.Clear()
